$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J45").Value = 0.2150495036779461
$ws.Range("I46").Value = 0.24
$ws.Range("H47").Value = 0.3087982760018804
$ws.Range("G48").Value = 0.32
$ws.Range("F49").Value = 0.4476495795507702
$ws.Range("E50").Value = 0.1088966743764388
$ws.Range("D51").Value = 0.1461563307127136
$ws.Range("C52").Value = 0.09547648014918764
$ws.Range("B53").Value = 0.0959495356205764
